# Update localization status report: files 25b24e84... and 8d543509...
# moved from "Ready for handoff" to "In Translation" status, while
# f6534ced... remains "Ready for handoff".

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold per-locale status ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "In Translation"   # 25b24e84...md, zh-cn
$wsOverview.Range("F3").Value = "In Translation"   # 25b24e84...md, de-de
$wsOverview.Range("E4").Value = "In Translation"   # 8d543509...md, zh-cn
$wsOverview.Range("F4").Value = "In Translation"   # 8d543509...md, de-de

# --- zh-cn sheet: column C holds Status ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"       # 25b24e84...md
$wsZhCn.Range("C4").Value = "In Translation"       # 8d543509...md

# --- de-de sheet: column C holds Status ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"       # 25b24e84...md
$wsDeDe.Range("C4").Value = "In Translation"       # 8d543509...md
